# Add a new worksheet "Emission_Factors" as the last sheet in the workbook
# and populate it with district -> emission factor data, matching the
# existing "header row" styling used throughout this workbook (bold,
# centered/top-aligned, thin border on all sides).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Emission_Factors"

# Header row
$ws.Range("A1").Value = "District"
$ws.Range("B1").Value = "Emission Factor"

$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# District -> emission factor data
$data = @(
    @("Chemnitz, Stadt", 0.00055956),
    @("Erzgebirgskreis", 0.00072968),
    @("Mittelsachsen", 0.00067033),
    @("Vogtlandkreis", 0.00049543),
    @("Zwickau", 0.00069221),
    @("Dresden, Stadt", 0.00126389),
    @("Bautzen", 0.00066184),
    @("Görlitz", 0.00055464),
    @("Meißen", 0.0005383),
    @("Sächsische Schweiz-Osterzgebirge", 0.00054913),
    @("Leipzig, Stadt", 0.00138366),
    @("Leipzig", 0.00058387),
    @("Nordsachsen", 0.00044573)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
